$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 2650
$ws.Range("I12").Value = 3366.6667
$ws.Range("K12").Value = 3366.6667
$ws.Range("M12").Value = -3196.6667

$ws.Range("H64").Value = 3115.8333
$ws.Range("J64").Value = 3157
$ws.Range("L64").Value = 3157
$ws.Range("N64").Value = -3653

$ws.Range("H67").Value = 3115.8333
$ws.Range("J67").Value = 3157
$ws.Range("L67").Value = 3157
$ws.Range("N67").Value = -4873

$ws.Range("H116").Value = 8417.5
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 8417.5
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 8417.5
$ws.Range("M116").Value = $null
$ws.Range("N116").Value = -15301.5

$ws.Range("H121").Value = 9178.571
$ws.Range("I121").Value = 700
$ws.Range("J121").Value = 9830.77
$ws.Range("K121").Value = 2100
$ws.Range("L121").Value = 29492.31
$ws.Range("M121").Value = -353
$ws.Range("N121").Value = -32986.31

$ws.Range("H137").Value = 1801.6428
$ws.Range("I137").Value = 1421.375
$ws.Range("K137").Value = 4264.125
$ws.Range("M137").Value = -1714.125

$ws.Range("H138").Value = 3050.0571
$ws.Range("I138").Value = 2862.125
$ws.Range("J138").Value = 3105.7407
$ws.Range("K138").Value = 8586.375
$ws.Range("L138").Value = 9317.222099999999
$ws.Range("M138").Value = -3446.375
$ws.Range("N138").Value = -19597.2221

$ws.Range("H141").Value = 2959.8572
$ws.Range("I141").Value = 2456.2942
$ws.Range("K141").Value = 7368.882599999999
$ws.Range("M141").Value = -2188.882599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1459.6666
$ws.Range("I61").Value = 1405.7222
$ws.Range("K61").Value = 1405.7222
$ws.Range("M61").Value = -1193.7222

$ws.Range("H74").Value = 40001570
$ws.Range("I74").Value = 83333920
$ws.Range("K74").Value = 83333920
$ws.Range("M74").Value = -83333046

$ws.Range("H77").Value = 40001570
$ws.Range("I77").Value = 83333920
$ws.Range("K77").Value = 416669600
$ws.Range("M77").Value = -416665232

$ws.Range("H97").Value = 100001176
$ws.Range("I97").Value = 1001.4286
$ws.Range("K97").Value = 1001.4286
$ws.Range("M97").Value = -505.4286

$ws.Range("H110").Value = 891.5
$ws.Range("I110").Value = 837.5
$ws.Range("J110").Value = 999.5
$ws.Range("K110").Value = 837.5
$ws.Range("L110").Value = 999.5
$ws.Range("M110").Value = 1207.5
$ws.Range("N110").Value = -5089.5

$ws.Range("H122").Value = 2112.353
$ws.Range("I122").Value = 1619.4375
$ws.Range("K122").Value = 4858.3125
$ws.Range("M122").Value = -2408.3125

$ws.Range("H132").Value = 15336.189
$ws.Range("I132").Value = 1600.4482
$ws.Range("K132").Value = 4801.3446
$ws.Range("M132").Value = -2271.3446

$ws.Range("H136").Value = 1459.6666
$ws.Range("I136").Value = 1405.7222
$ws.Range("K136").Value = 4217.1666
$ws.Range("M136").Value = -1667.1666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").Value = $null

$ws.Range("H94").Value = 994.5
$ws.Range("I94").Value = 471.26666
$ws.Range("K94").Value = 471.26666
$ws.Range("M94").Value = -20.26666

$ws.Range("H134").Value = 2626.4314
$ws.Range("I134").Value = 2725.8262
$ws.Range("K134").Value = 8177.4786
$ws.Range("M134").Value = -5642.4786

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 1008
$ws.Range("I17").Value = 1008
$ws.Range("K17").Value = 1008
$ws.Range("M17").Value = -834

$ws.Range("H19").Value = 173
$ws.Range("I19").Value = 173
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 173
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -3
$ws.Range("N19").Value = $null

$ws.Range("H24").Value = 173
$ws.Range("I24").Value = 173
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 173
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -3
$ws.Range("N24").Value = $null

$ws.Range("H108").Value = 8584
$ws.Range("I108").Value = 8584
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 8584
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = -4744
$ws.Range("N108").Value = $null

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").Value = $null

$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = $null

$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").Value = $null

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = $null

$ws.Range("H121").Value = 5800
$ws.Range("I121").Value = 5800
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 5800
$ws.Range("L121").Value = 0
$ws.Range("M121").Value = -4490
$ws.Range("N121").Value = $null

$ws.Range("H122").Value = 1531.2
$ws.Range("I122").Value = 1671.25
$ws.Range("J122").Value = 1371.1428
$ws.Range("K122").Value = 5013.75
$ws.Range("L122").Value = 4113.428400000001
$ws.Range("M122").Value = -2563.75
$ws.Range("N122").Value = -9013.428400000001

$ws.Range("H132").Value = 3033.76
$ws.Range("I132").Value = 2282.4375
$ws.Range("J132").Value = 4369.4443
$ws.Range("K132").Value = 6847.3125
$ws.Range("L132").Value = 13108.3329
$ws.Range("M132").Value = -4317.3125
$ws.Range("N132").Value = -18168.3329

$ws.Range("H133").Value = 38845
$ws.Range("J133").Value = 38845
$ws.Range("L133").Value = 38845
$ws.Range("N133").Value = -43905

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 57.1
$ws.Range("J12").Value = 78.71429000000001
$ws.Range("L12").Value = 236.14287
$ws.Range("N12").Value = -582.14287

$ws.Range("H33").Value = 60.857143
$ws.Range("J33").Value = 67.666664
$ws.Range("L33").Value = 405.999984
$ws.Range("N33").Value = -971.999984

$ws.Range("H131").Value = 715.9
$ws.Range("J131").Value = 720
$ws.Range("L131").Value = 2160
$ws.Range("N131").Value = -12240

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 16672500
$ws.Range("J52").Value = 16672500
$ws.Range("L52").Value = 16672500
$ws.Range("N52").Value = -16673018

$ws.Range("H80").Value = 4466.231
$ws.Range("I80").Value = 3281
$ws.Range("K80").Value = 3281
$ws.Range("M80").Value = -2283

$ws.Range("H83").Value = 4466.231
$ws.Range("I83").Value = 3281
$ws.Range("K83").Value = 16405
$ws.Range("M83").Value = -11413

$ws.Range("H97").Value = 2359.3845
$ws.Range("I97").Value = 2390.25
$ws.Range("J97").Value = 2310
$ws.Range("K97").Value = 2390.25
$ws.Range("L97").Value = 2310
$ws.Range("M97").Value = -1894.25
$ws.Range("N97").Value = -3302

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 8800.333000000001
$ws.Range("I32").Value = 8800.333000000001
$ws.Range("K32").Value = 8800.333000000001
$ws.Range("M32").Value = -8483.333000000001

$ws.Range("H82").Value = 2700
$ws.Range("I82").Value = 3266.6667
$ws.Range("J82").Value = 1000
$ws.Range("K82").Value = 3266.6667
$ws.Range("L82").Value = 1000
$ws.Range("M82").Value = -2905.6667
$ws.Range("N82").Value = -1722

$ws.Range("H85").Value = 2700
$ws.Range("I85").Value = 3266.6667
$ws.Range("J85").Value = 1000
$ws.Range("K85").Value = 3266.6667
$ws.Range("L85").Value = 1000
$ws.Range("M85").Value = -2018.6667
$ws.Range("N85").Value = -3496

$ws.Range("H93").Value = 1549.5
$ws.Range("I93").Value = 874.25
$ws.Range("J93").Value = 2900
$ws.Range("K93").Value = 874.25
$ws.Range("L93").Value = 2900
$ws.Range("M93").Value = 373.75
$ws.Range("N93").Value = -5396

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1258.1
$ws.Range("I107").Value = 713.5
$ws.Range("K107").Value = 2140.5
$ws.Range("M107").Value = -220.5

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = $null

$ws.Range("H110").Value = 35000
$ws.Range("J110").Value = 35000
$ws.Range("L110").Value = 35000
$ws.Range("N110").Value = -43180

$ws.Range("H122").Value = 899.0769
$ws.Range("I122").Value = 862.375
$ws.Range("K122").Value = 2587.125
$ws.Range("M122").Value = -137.125

$ws.Range("H136").Value = 19609914
$ws.Range("I136").Value = 27778836
$ws.Range("J136").Value = 4499.1333
$ws.Range("K136").Value = 83336508
$ws.Range("L136").Value = 13497.3999
$ws.Range("M136").Value = -83333958
$ws.Range("N136").Value = -18597.3999
